$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.258.88"
$ws.Range("E2").Value = "  +4.56%  "
$ws.Range("D3").Value = "2.729.91"
$ws.Range("E3").Value = "  +4.07%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "530.19"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.84%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "147.88"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +2.15%  "
$ws.Range("E7").Value = "  -0.20%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.582"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.59%  "
$ws.Range("D9").Value = "2.759.17"
$ws.Range("E9").Value = "  +5.08%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "7.18"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +14.79%  "
$ws.Range("E11").Value = "  +2.33%  "
$ws.Range("E12").Value = "  +2.67%  "
$ws.Range("E13").Value = "  +3.22%  "
$ws.Range("D14").Value = "3.203.45"
$ws.Range("E14").Value = "  +4.07%  "
$ws.Range("D15").Value = "61.197.36"
$ws.Range("E15").Value = "  +4.41%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "21.62"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +4.50%  "
$ws.Range("D17").Value = "2.800.62"
$ws.Range("E17").Value = "  +6.73%  "
$ws.Range("E18").Value = "  +2.38%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "347.53"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.70%  "
$ws.Range("E20").Value = "  +2.42%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "10.67"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +5.03%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.48"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +6.02%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.996"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.36%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "63.75"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +3.89%  "
$ws.Range("E25").Value = "  +5.61%  "
$ws.Range("E26").Value = "  +2.25%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.994"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.44%  "
$ws.Range("D28").Value = "0.0₃0832"
$ws.Range("E28").Value = "  +4.22%  "
$ws.Range("E29").Value = "  +5.64%  "
$ws.Range("E30").Value = "  +10.03%  "
$ws.Range("E31").Value = "  -0.20%  "
$ws.Range("E32").Value = "  +2.43%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "19.17"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +1.85%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "150.20"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.26%  "
$ws.Range("E35").Value = "  +8.40%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.24"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +9.74%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.926"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -4.85%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.908"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +9.16%  "
$ws.Range("E39").Value = "  +9.25%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "37.27"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +1.94%  "
$ws.Range("E41").Value = "  +2.66%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "285.59"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +2.54%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "20.58"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +5.72%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.626"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +5.07%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0994"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +1.30%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.12%  "
$ws.Range("D47").Value = "2.137.30"
$ws.Range("E47").Value = "  +8.08%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "4.98"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +7.89%  "
$ws.Range("E49").Value = "  +5.19%  "
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "10.54"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +2.27%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "19.50"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +6.83%  "
